# chore: adapt column header formatting to respective input file names
#
# Renames the `_old` / `_new` header-name suffixes to `_FV2404` / `_FV2410`
# (the two EDI@Energy "Formatversion" identifiers being diffed), wraps the
# sheet's data range in a real Excel Table ("Table1"), and freezes the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rewrite the header row (A1:U1) -------------------------------------
# Column order on the sheet is: 10 "_old" columns, a "diff" column, then
# 10 "_new" columns -> map 1:1 onto "_FV2404" / "_FV2410".
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into an Excel Table -----------------------------
$usedRange = $ws.Range("A1:U71")
$table = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("A1").Select()
